# Apply cryptos.xlsx update: Tue Apr 16 06:43:08 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.937.69"
$ws.Range("E2").Value = "  -4.99%  "
$ws.Range("D3").Value = "3.070.11"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.68%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.063.22"
$ws.Range("E8").Value = "  -5.24%  "
$ws.Range("E9").Value = "  -4.90%  "
$ws.Range("E10").Value = "  -7.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.05%  "
$ws.Range("E12").Value = "  -5.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.97%  "
$ws.Range("E14").Value = "  -9.11%  "
$ws.Range("D15").Value = "3.563.52"
$ws.Range("E15").Value = "  -5.24%  "
$ws.Range("D16").Value = "62.971.02"
$ws.Range("E16").Value = "  -4.98%  "
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "3.069.12"
$ws.Range("E18").Value = "  -5.38%  "
$ws.Range("E19").Value = "  -5.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "486.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.49%  "
$ws.Range("E22").Value = "  -4.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.74%  "
$ws.Range("E28").Value = "  -5.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -14.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.23%  "
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "518.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.80%  "
$ws.Range("E37").Value = "  -12.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0402"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -14.71%  "
$ws.Range("D39").Value = "3.100.41"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.97%  "
$ws.Range("E41").Value = "  -5.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.75%  "
$ws.Range("E44").Value = "  -8.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.70%  "
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0497"
$ws.Range("E50").Value = "  -13.64%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +49.19%  "
